# Add account-information rows (module "test2".."test42") beneath the existing
# "test1" property row, mirroring the Id/Type/Public/.../Level/test1 block above it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A uses the same style as the existing "test1" row (A11) for every new row.
$ws.Range("A11").Copy()
$ws.Range("A12:A52").PasteSpecial(-4122)

# Most column-B cells are plain numbers (style carried over from the column default),
# but a handful of rows are "module header" rows that reuse B11's text style ("1" as a
# shared string) instead of the plain numeric style.
$ws.Range("B11").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("B31").PasteSpecial(-4122)
$ws.Range("B33").PasteSpecial(-4122)
$ws.Range("B38").PasteSpecial(-4122)
$ws.Range("B43").PasteSpecial(-4122)
$ws.Range("B45").PasteSpecial(-4122)
$ws.Range("B50").PasteSpecial(-4122)

# Populate column A (row labels) for the new rows 12-52.
$ws.Range("A12").Value = "test2"
$ws.Range("A13").Value = "test3"
$ws.Range("A14").Value = "test4"
$ws.Range("A15").Value = "test5"
$ws.Range("A16").Value = "test6"
$ws.Range("A17").Value = "test7"
$ws.Range("A18").Value = "test8"
$ws.Range("A19").Value = "test9"
$ws.Range("A20").Value = "test10"
$ws.Range("A21").Value = "test11"
$ws.Range("A22").Value = "test12"
$ws.Range("A23").Value = "test13"
$ws.Range("A24").Value = "test14"
$ws.Range("A25").Value = "test15"
$ws.Range("A26").Value = "test16"
$ws.Range("A27").Value = "test17"
$ws.Range("A28").Value = "test18"
$ws.Range("A29").Value = "test19"
$ws.Range("A30").Value = "test20"
$ws.Range("A31").Value = "test21"
$ws.Range("A32").Value = "test22"
$ws.Range("A33").Value = "test23"
$ws.Range("A34").Value = "test24"
$ws.Range("A35").Value = "test25"
$ws.Range("A36").Value = "test26"
$ws.Range("A37").Value = "test27"
$ws.Range("A38").Value = "test28"
$ws.Range("A39").Value = "test29"
$ws.Range("A40").Value = "test30"
$ws.Range("A41").Value = "test31"
$ws.Range("A42").Value = "test32"
$ws.Range("A43").Value = "test33"
$ws.Range("A44").Value = "test34"
$ws.Range("A45").Value = "test35"
$ws.Range("A46").Value = "test36"
$ws.Range("A47").Value = "test37"
$ws.Range("A48").Value = "test38"
$ws.Range("A49").Value = "test39"
$ws.Range("A50").Value = "test40"
$ws.Range("A51").Value = "test41"
$ws.Range("A52").Value = "test42"

# Populate column B: "1" (text) for module-header rows, 1 (number) for the rest.
$ws.Range("B16").Value = "1"
$ws.Range("B21").Value = "1"
$ws.Range("B26").Value = "1"
$ws.Range("B31").Value = "1"
$ws.Range("B33").Value = "1"
$ws.Range("B38").Value = "1"
$ws.Range("B43").Value = "1"
$ws.Range("B45").Value = "1"
$ws.Range("B50").Value = "1"

$ws.Range("B12").Value = 1
$ws.Range("B13").Value = 1
$ws.Range("B14").Value = 1
$ws.Range("B15").Value = 1
$ws.Range("B17").Value = 1
$ws.Range("B18").Value = 1
$ws.Range("B19").Value = 1
$ws.Range("B20").Value = 1
$ws.Range("B22").Value = 1
$ws.Range("B23").Value = 1
$ws.Range("B24").Value = 1
$ws.Range("B25").Value = 1
$ws.Range("B27").Value = 1
$ws.Range("B28").Value = 1
$ws.Range("B29").Value = 1
$ws.Range("B30").Value = 1
$ws.Range("B32").Value = 1
$ws.Range("B34").Value = 1
$ws.Range("B35").Value = 1
$ws.Range("B36").Value = 1
$ws.Range("B37").Value = 1
$ws.Range("B39").Value = 1
$ws.Range("B40").Value = 1
$ws.Range("B41").Value = 1
$ws.Range("B42").Value = 1
$ws.Range("B44").Value = 1
$ws.Range("B46").Value = 1
$ws.Range("B47").Value = 1
$ws.Range("B48").Value = 1
$ws.Range("B49").Value = 1
$ws.Range("B51").Value = 1
$ws.Range("B52").Value = 1

# Match the author's final selection (cell B33) recorded in the saved sheet view.
[void]$ws.Range("B33").Select()
